$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with data, mirroring the formatting of row 4 (date style in col A, inline strings elsewhere)
$ws.Range("A5").Value2 = $ws.Range("A4").Value2
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Range("B5").Value = "Thayna Silva Santana"
$ws.Range("C5").Value = "Aristides Attico 2"
$ws.Range("D5").Value = "thayna.santana"
$ws.Range("E5").Value = "Francisco"
$ws.Range("F5").Value = "Estagiaria"
$ws.Range("G5").Value = "Automação"
$ws.Range("H5").Value = "Todos"
$ws.Range("I5").Value = "Notebook + Carregador"
$ws.Range("J5").Value = "LUM-001-001-078"
$ws.Range("K5").Value = "Mayara Almeida"
$ws.Range("L5").Value = "23/05/2024 17:46:56"

$wb.Save()
